$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column E (Statistical analysis), shifting old E..F to F..G ---
$ws.Columns("E:E").Insert()

# --- Row 1 (header) ---
$ws.Range("E1").Value = "Statistical analysis"

# --- Row 2 (Coufal / Silene acaulis) ---
$ws.Range("E2").Value = "presence/pseudo-absence GLM "

# --- Row 3 (Fer / Vipera berus) ---
$ws.Range("D3").Value = "Are the predictors good for the distribution of V. berus?"
$ws.Range("D3").Characters(49, 8).Font.Italic = $true
$ws.Range("D3").Characters(49, 8).Font.Name = "Calibri"
$ws.Range("D3").Characters(49, 8).Font.Size = 11
$ws.Range("D3").Characters(57, 1).Font.Name = "Calibri"
$ws.Range("D3").Characters(57, 1).Font.Size = 11
$ws.Range("E3").Value = "presence/pseudo-absence GLM "

# --- Row 4 (Gulkova / Ursus arctus) ---
$ws.Range("C4").Font.Italic = $true
$ws.Range("E4").Value = "MaxEnt"
$ws.Range("F4").Value = "Europe"

# --- Row 5 (Kucera / Spatially explicit neutral model simulation) ---
$ws.Range("D5").Value = "How does # of species develop through time, and how does it depend on some simulation parameters?"
$ws.Range("E5").Value = "Graphical illustration"

# --- Row 7 (Sejk / The bald eagle) ---
$ws.Range("E7").Value = "abundance Poisson GLM"

# --- Row 8 (Stepanova / European beaver) ---
$ws.Range("E8").Value = "MaxEnt "

# --- Row 9 (Valek / new project: Species richness of mammals) ---
$ws.Range("C9").Value = "Species richness of mammals"
$ws.Range("C9").Font.Italic = $true
$ws.Range("D9").Value = "How does # of species develop through time, and how does it depend on some simulation parameters?"
$ws.Range("E9").Value = "linear or poisson regression"
$ws.Range("F9").Value = "Czech Republic"
$ws.Range("G9").Value = "time"

# --- Row 11 (Voska / Limulus polyphaemus) ---
$ws.Range("E11").Value = "MaxEnt"

# --- Row 12 (Zak / Lama guanicoe) ---
$ws.Range("E12").Value = "MaxEnt"

# --- Column widths ---
$ws.Range("D1").ColumnWidth = 85.269
$ws.Range("E1").ColumnWidth = 34.432
$ws.Range("F1").ColumnWidth = 25.432
$ws.Range("G1").ColumnWidth = 56.595
$ws.Range("H1").ColumnWidth = 24.269

# --- Selection ---
$ws.Range("D10").Select()
